$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A blank-but-formatted cell shows up at C1 (same default style as B1/E1) after
# the edit -- reproduce it by cloning B1's (empty, default-styled) format onto C1.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Shipping cost for AMS1117 (row 7) gets allocated 10 in the "Tung" column (G).
# G9/E10 totals are formulas and recalculate automatically from this.
$ws.Range("G7").Value = 10

# Reflect the new active selection recorded in the sheet view.
$ws.Range("H13").Select()
